$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MATLAB stress value (B7) with the more rigorous value from MATLAB
$ws.Range("B7").Value = 72

# Move active selection/cursor to B8 (reflects where the user ended up editing)
$ws.Range("B8").Select()
